$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.388.84'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '1.560.57'
$ws.Range("E3").Value = '  -0.80%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '285.74'
$ws.Range("E6").Value = '  -1.84%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3637'
$ws.Range("E7").Value = '  -2.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.68'
$ws.Range("E8").Value = '  -2.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3333'
$ws.Range("E9").Value = '  -2.12%  '
$ws.Range("E10").Value = '  -1.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07381'
$ws.Range("E11").Value = '  -2.59%  '
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.74'
$ws.Range("E13").Value = '  -3.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.909'
$ws.Range("E14").Value = '  -1.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.856'
$ws.Range("E15").Value = '  -1.45%  '
$ws.Range("D16").Value = '1.560.63'
$ws.Range("E17").Value = '  -2.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '88.64'
$ws.Range("E18").Value = '  -2.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06726'
$ws.Range("E19").Value = '  -0.22%  '
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.311'
$ws.Range("E21").Value = '  +0.56%  '
$ws.Range("E22").Value = '  -2.72%  '
$ws.Range("E23").Value = '  -2.60%  '
$ws.Range("D24").Value = '22.374.44'
$ws.Range("E24").Value = '  -0.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.387'
$ws.Range("E25").Value = '  +2.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.545'
$ws.Range("E26").Value = '  -1.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '149.09'
$ws.Range("E27").Value = '  +0.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.35'
$ws.Range("E28").Value = '  -4.06%  '
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.71'
$ws.Range("E30").Value = '  -2.61%  '
$ws.Range("D31").Value = '1.735.71'
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.051'
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.090'
$ws.Range("E33").Value = '  -0.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.987'
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.554'
$ws.Range("E35").Value = '  -3.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08222'
$ws.Range("E36").Value = '  -2.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02366'
$ws.Range("E37").Value = '  -4.06%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.294'
$ws.Range("E38").Value = '  -6.63%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2205'
$ws.Range("E39").Value = '  -4.09%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06336'
$ws.Range("E40").Value = '  -3.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.301'
$ws.Range("E41").Value = '  -3.47%  '
$ws.Range("E42").Value = '  -2.75%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6042'
$ws.Range("E43").Value = '  -3.86%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.54'
$ws.Range("E45").Value = '  -3.38%  '
$ws.Range("E46").Value = '  -1.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5714'
$ws.Range("E47").Value = '  -2.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.24'
$ws.Range("E48").Value = '  -4.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.994'
$ws.Range("E49").Value = '  -4.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.204'
$ws.Range("E50").Value = '  -2.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07212'
$ws.Range("E51").Value = '  -1.59%  '
